$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (Kavish -> Sample1)
$ws.Range("A3").Value = "Sample1"
$ws.Range("B3").Value = "sample1"
$ws.Range("C3").Value = "sample1234"

# Update row 4 (Preet -> Sample2), reuse "sample1234" as the password value
$ws.Range("A4").Value = "Sample2"
$ws.Range("B4").Value = "sample2"
$ws.Range("C4").Value = "sample1234"

# Update the active selection to D4
$ws.Range("D4").Select()
